$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the lingering N4 selection back to the default A1 (the diff drops the
# explicit <selection> from sheetView; this is the closest equivalent the
# object model exposes).
$ws.Range("A1").Select() | Out-Null

# --- Column A:C width: 36.7109375 -> 37.140625 (closest reachable grid value) ---
$ws.Columns("A:C").ColumnWidth = 36.25

# --- Row heights ---
$ws.Rows(1).RowHeight = 55.5
$ws.Rows(2).RowHeight = 13.5
$ws.Rows(3).RowHeight = 13.5

# --- New column N data (extends table with 2023 figures) ---

# Row 3: empty divider cell matching M3's style (border/thick-bottom row)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# Row 4: year header 2023, matching style of M4 (year header cell)
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

# Row 5: data value (style already matches default column style s=1)
$ws.Range("N5").Value = 4.3499999999999996

# Row 6: data value (style already matches default column style s=1)
$ws.Range("N6").Value = 4.3499999999999996

# Row 7: "-" placeholder, right aligned (new style: font 1, horizontal=right, vertical=center)
$ws.Range("N7").HorizontalAlignment = -4152
$ws.Range("N7").Value = "-"

# Row 8: "-" placeholder matching M8's style (bottom border + right aligned)
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = "-"
